$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7666
$ws.Range("J3").Value = 8047
$ws.Range("J4").Value = 1749
$ws.Range("J5").Value = 624
$ws.Range("J6").Value = 11007
$ws.Range("J7").Value = 29093

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 232
$ws.Range("J4").Value = 132
$ws.Range("J5").Value = 86
$ws.Range("J6").Value = 223
$ws.Range("J7").Value = 829
$ws.Range("J8").Value = 1842
$ws.Range("J9").Value = 144
$ws.Range("J11").Value = 531
$ws.Range("J14").Value = 153
$ws.Range("J20").Value = 629
$ws.Range("J23").Value = 267
$ws.Range("J25").Value = 153
$ws.Range("J27").Value = 178
$ws.Range("J29").Value = 1554
$ws.Range("J33").Value = 1310
$ws.Range("J34").Value = 134
$ws.Range("J35").Value = 35
$ws.Range("J36").Value = 397
$ws.Range("J37").Value = 900
$ws.Range("J40").Value = 66
$ws.Range("J42").Value = 1235
$ws.Range("J44").Value = 232
$ws.Range("J48").Value = 319
$ws.Range("J49").Value = 175
$ws.Range("J51").Value = 365
$ws.Range("J53").Value = 437
$ws.Range("J54").Value = 568
$ws.Range("J55").Value = 455
$ws.Range("J59").Value = 36
$ws.Range("J60").Value = 173
$ws.Range("J63").Value = 84
$ws.Range("J64").Value = 191
$ws.Range("J67").Value = 1057
$ws.Range("J76").Value = 412
$ws.Range("J79").Value = 797
$ws.Range("J84").Value = 240
$ws.Range("J85").Value = 1188
$ws.Range("J86").Value = 175
$ws.Range("J87").Value = 97
$ws.Range("J88").Value = 311
$ws.Range("J89").Value = 365
$ws.Range("J90").Value = 306
$ws.Range("J91").Value = 334
$ws.Range("J93").Value = 122
$ws.Range("J94").Value = 326
$ws.Range("J95").Value = 416
$ws.Range("J98").Value = 214
$ws.Range("J99").Value = 443
$ws.Range("J100").Value = 50
$ws.Range("J101").Value = 29093

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J2").Value = 50
$ws.Range("J3").Value = 30
$ws.Range("J7").Value = 153

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 264
$ws.Range("J7").Value = 829

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 257
$ws.Range("J7").Value = 531

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J4").Value = 37
$ws.Range("J7").Value = 365

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 316
$ws.Range("J7").Value = 1188

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 78
$ws.Range("J7").Value = 437

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 485
$ws.Range("J3").Value = 527
$ws.Range("J5").Value = 48
$ws.Range("J6").Value = 687
$ws.Range("J7").Value = 1842

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 295
$ws.Range("J3").Value = 434
$ws.Range("J6").Value = 470
$ws.Range("J7").Value = 1310

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 150
$ws.Range("J6").Value = 87
$ws.Range("J7").Value = 416

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 272
$ws.Range("J3").Value = 303
$ws.Range("J6").Value = 261
$ws.Range("J7").Value = 900

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 122
$ws.Range("J7").Value = 443

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 398
$ws.Range("J6").Value = 292
$ws.Range("J7").Value = 1057

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J6").Value = 78
$ws.Range("J7").Value = 240

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J4").Value = 46
$ws.Range("J6").Value = 261
$ws.Range("J7").Value = 568

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 470
$ws.Range("J3").Value = 547
$ws.Range("J4").Value = 84
$ws.Range("J6").Value = 396
$ws.Range("J7").Value = 1554

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J3").Value = 61
$ws.Range("J7").Value = 319

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J2").Value = 72
$ws.Range("J7").Value = 232

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 75
$ws.Range("J3").Value = 91
$ws.Range("J6").Value = 212
$ws.Range("J7").Value = 412

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J2").Value = 66
$ws.Range("J3").Value = 54
$ws.Range("J7").Value = 223

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 250
$ws.Range("J6").Value = 654
$ws.Range("J7").Value = 1235

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J2").Value = 87
$ws.Range("J6").Value = 257
$ws.Range("J7").Value = 455

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 267

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 92
$ws.Range("J7").Value = 334

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 232
$ws.Range("J7").Value = 797

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J3").Value = 47
$ws.Range("J7").Value = 191

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 173
$ws.Range("J7").Value = 629

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 131
$ws.Range("J7").Value = 397

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 122

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 50

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J2").Value = 61
$ws.Range("J6").Value = 175
$ws.Range("J7").Value = 326

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 66
$ws.Range("J7").Value = 153

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J6").Value = 134
$ws.Range("J7").Value = 214

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J3").Value = 48
$ws.Range("J7").Value = 144

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J6").Value = 88
$ws.Range("J7").Value = 232

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J3").Value = 68
$ws.Range("J6").Value = 164
$ws.Range("J7").Value = 311

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J3").Value = 29
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 110
$ws.Range("J7").Value = 306

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J2").Value = 78
$ws.Range("J7").Value = 365

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 173

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 66

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 97
